$d = $word.ActiveDocument

# The original paragraph is a single run:
#   "Aaaa bbbb cccc dddd eeee ffff gggg hhhh iiii jjjj kkkk llll mmmm nnnn oooo="
# It must become six separate runs (same Consolas formatting) whose text
# concatenates to "Map counts nearby unvisited and unphotographed cities."
#
# A plain Find/Replace (or sequential Range.InsertAfter calls) would leave a
# single merged run, so instead we replace the paragraph's content with raw
# OOXML via Range.InsertXML, which preserves each <w:r> as its own run.

$para = $d.Paragraphs.Item(1)
$r = $para.Range
# Exclude the trailing paragraph mark from the range so only the run content
# is replaced.
$r.MoveEnd(1, -1) | Out-Null

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t xml:space="preserve">Map </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t>counts</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t xml:space="preserve">nearby </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t>unvisited and unphotographed cities</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t>.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$r.InsertXML($xml)
